$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto push: a fresh data point for 2026/01/12 was appended to the
# rolling log. New rows are always inserted just above the first
# "2026/12/29" block (row 619), pushing the rest of the table down by one
# row (old row 660 -> new row 661), which is why <dimension> grows from
# A1:D660 to A1:D661.
$ws.Rows("619:619").Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real
# date serial. A bare .Value = "2026/01/12" would get auto-parsed into a
# date by Excel's smart-entry logic, so prefix with an apostrophe to force
# text entry, exactly like typing it in by hand.
$ws.Range("A619").Value = "'2026/01/12"
$ws.Range("B619").Value = "月"
$ws.Range("C619").Value = 7
$ws.Range("D619").Value = 201
